$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the index column (column A held the numeric row index 0..9);
# deleting it shifts all remaining columns (the "people" header, weekday
# headers, date headers, staff names and their data) one column to the left,
# carrying each cell's existing formatting/style along with it.
$ws.Columns.Item(1).Delete()

# The original data had a corrupted date in what is now column V of the
# header row (it read as 2/27/2020 instead of 2/27/2021, breaking the
# otherwise-sequential weekly date headers). Correct it now that the index
# column has been removed, using the Excel serial date value for 2/27/2021.
$ws.Range("V1").Value = 44254
